$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 124.8841655586322
$ws.Range("C3").Value = 124.1905198343063
$ws.Range("C4").Value = 201.0397952678128
$ws.Range("C5").Value = 209.6366870621295
$ws.Range("C6").Value = 415.5888154020474
$ws.Range("C7").Value = 452.9531410628954
$ws.Range("C8").Value = 313.1212541168178
$ws.Range("C9").Value = 398.4469012570319
$ws.Range("C10").Value = 365.267592606145
$ws.Range("C11").Value = 354.8472065463146
$ws.Range("C12").Value = 541.0434353570234
$ws.Range("C13").Value = 503.3169421585232
$ws.Range("C14").Value = 919.9346972302124
$ws.Range("C15").Value = 673.1850647811689
$ws.Range("C16").Value = 516.3256785400293
$ws.Range("C17").Value = 563.25236871781
$ws.Range("C18").Value = 981.1130843061849
$ws.Range("C19").Value = 1514.443225538084
$ws.Range("C20").Value = 813.6677088815854
$ws.Range("C21").Value = 756.7779807417407
$ws.Range("C22").Value = 781.9314304310158
$ws.Range("C23").Value = 796.3158327344975
$ws.Range("C24").Value = 1097.429221319476
$ws.Range("C25").Value = 928.853209725026
$ws.Range("C26").Value = 928.2380526109997
$ws.Range("C27").Value = 1105.619223133979
$ws.Range("C28").Value = 986.5509504907519
$ws.Range("C29").Value = 899.3294493845233
$ws.Range("C30").Value = 1363.829707685258
$ws.Range("C31").Value = 1106.411639163403
$ws.Range("C32").Value = 975.7971719791447
$ws.Range("C33").Value = 1084.610957170714
$ws.Range("C34").Value = 1231.70890405144
$ws.Range("C35").Value = 1238.981715426326
$ws.Range("C36").Value = 1130.765106308259
$ws.Range("C37").Value = 1851.066220260111
$ws.Range("C38").Value = 2678.07312783917
$ws.Range("C39").Value = 1695.5922174267
$ws.Range("C40").Value = 1848.241111526819
$ws.Range("C41").Value = 1563.604902509132
$ws.Range("C42").Value = 1627.201143406835
$ws.Range("C43").Value = 1566.49703229932
$ws.Range("C44").Value = 1860.480274610089
$ws.Range("C45").Value = 1653.434195468234
$ws.Range("C46").Value = 2277.482947079543
$ws.Range("C47").Value = 2277.428982068564
$ws.Range("C48").Value = 1771.536879235642
$ws.Range("C49").Value = 1504.431802963491
$ws.Range("C50").Value = 3386.155625376966
